$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Requisitos:"

$reqText = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("B23").Value = $reqText
$ws.Range("C23").Value = $reqText

$ws.Rows.Item(23).RowHeight = 30

$ws.Range("B21").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
